# Add a new cdtede command: "set_can_broadcast_readout_stop"
# Inserted as a new row 17 on the "all_systems" sheet, pushing the
# existing rows 17-23 down to 18-24.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 17, copying row 17's current formatting so the
# inserted row inherits identical styles/number-formats.
$ws.Rows(17).Copy()
$ws.Rows(17).Insert()

# The Insert() above only partially preserves per-cell styles, so
# reinforce the formatting by pasting the (now shifted-down) original
# row 17 -- which now lives at row 18 -- format over the new row 17.
$ws.Range("A18:AC18").Copy()
$ws.Range("A17:AC17").PasteSpecial(-4122)  # xlPasteFormats

# Fill in the new command's data.
$ws.Range("A17").Value2 = "set_can_broadcast_readout_stop"
$ws.Range("D17").Value2 = 10111
$ws.Range("E17").Formula = '=_xlfn.CONCAT("0x", DEC2HEX(_xlfn.BITLSHIFT($C17,7) + BIN2DEC($D17)))'
$ws.Range("AA17").Value2 = "0x3c3c0101050505053c3c3c3c"
$ws.Range("AC17").Value2 = "ALL canisters stop reading out."

# Reflect the author's final cursor position/selection in the sheet view.
$ws.Activate()
$ws.Range("AC17").Select()
